# Codebook Tutorial edits
$d = $word.ActiveDocument

function Get-PkgXml($bodyXml) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>$bodyXml</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@
}

function Split-Run($needle, $bodyXml) {
    # Locate the run's exact text, delete it, then insert the replacement
    # run(s) at the now-collapsed point. Doing a plain InsertXML across a
    # range whose End sits exactly on a following sibling run/hyperlink can
    # re-order content, so delete-then-insert-at-collapsed-point is used
    # instead for reliable in-place splitting.
    $rng = $d.Content
    $rng.Find.Execute($needle) | Out-Null
    if (-not $rng.Find.Found) {
        throw "Not found: $needle"
    }
    $start = $rng.Start
    $target = $d.Range($start, $rng.End)
    $target.Delete()
    $collapsed = $d.Range($start, $start)
    $collapsed.InsertXML((Get-PkgXml $bodyXml))
}

$TNR = 'w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"'

# ---------------------------------------------------------------------------
# 1) Title: "Codebook for Dummies" -> "Codebook for" | bookmark _GoBack | " Dummies"
#    (the _GoBack bookmark is moved here from the Step 4 paragraph; adding it
#    by name automatically removes the old one since bookmark names are
#    unique within a document)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Codebook for Dummies") | Out-Null
$splitAt = $rng.Start + 12   # length of "Codebook for"
$insertPoint = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("_GoBack", $insertPoint) | Out-Null

# ---------------------------------------------------------------------------
# 2) Step 1: "Step 1: Go to Codebook website using " ->
#    "Step 1: Go to" | " the" | " Codebook website using "
# ---------------------------------------------------------------------------
$body = "<w:p><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>Step 1: Go to</w:t></w:r><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`"> the</w:t></w:r><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`"> Codebook website using </w:t></w:r></w:p>"
Split-Run "Step 1: Go to Codebook website using " $body

# ---------------------------------------------------------------------------
# 3) Step 2: `Upload data file on your computer using the "Browse" button at
#    the top of the page.` ->
#    "Upload " | "the " | "data file " | "from" |
#    ` your computer using the "Browse" button at the top of the page.`
# ---------------------------------------------------------------------------
$quoteOpen = [char]0x201C
$quoteClose = [char]0x201D
$needle2 = "Upload data file on your computer using the " + $quoteOpen + "Browse" + $quoteClose + " button at the top of the page."
$body = "<w:p><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`">Upload </w:t></w:r><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`">the </w:t></w:r><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`">data file </w:t></w:r><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>from</w:t></w:r><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`"> your computer using the $quoteOpen" + "Browse$quoteClose button at the top of the page.</w:t></w:r></w:p>"
Split-Run $needle2 $body

# ---------------------------------------------------------------------------
# 4) Step 3: the 3-space run right after "Generate codebook!" becomes a
#    closing curly quote (simple in-place text swap within the same run).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Step 3: Click") | Out-Null
$after = $d.Range($rng.End, $d.Content.End)
$after.Find.Execute("   ", $true, $false, $false, $false, $false, $true, 1, $false, [string]$quoteClose, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Step 4: split the sentence:
#    "Step 4: Make any changes on the created document by editing in the code
#     on the left side of the screen and clicking "Generate codebook!" again." ->
#    "Step 4: Make changes " | "to" | " the created document by editing the
#     code on the left side of the screen and clicking "Generate codebook!"
#     again."
#    (The old _GoBack bookmark that used to sit in this paragraph was already
#    relocated to the title above.)
# ---------------------------------------------------------------------------
$needle4 = "Step 4: Make any changes on the created document by editing in the code on the left side of the screen and clicking " + $quoteOpen + "Generate codebook!" + $quoteClose + " again."
$body = "<w:p><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`">Step 4: Make changes </w:t></w:r><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>to</w:t></w:r><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`"> the created document by editing the code on the left side of the screen and clicking $quoteOpen" + "Generate codebook!$quoteClose again.</w:t></w:r></w:p>"
Split-Run $needle4 $body

# ---------------------------------------------------------------------------
# 6) Step 5: append a new run containing "." after the existing sentence
# ---------------------------------------------------------------------------
$needle5 = "Step 5: When ready, click " + $quoteOpen + "Download codebook" + $quoteClose + " to save to your computer"
$body = "<w:p><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`">Step 5: When ready, click $quoteOpen" + "Download codebook$quoteClose to save to your computer</w:t></w:r><w:r><w:rPr><$TNR/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>.</w:t></w:r></w:p>"
Split-Run $needle5 $body

$d.Save()
